# Updated with Field Name with standard name and Minor release
#
# Rename the header cells in row 1 to the standard/short field names and
# normalize the HTTP status code values, matching the "standard name"
# relabeling described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1: RequestProcessingType -> RequestHeaders
$ws.Range("F1").Value = "RequestHeaders"

# I1: HTTPAction -> Action
$ws.Range("I1").Value = "Action"

# J1: ExcludeField -> ExcludeFields
$ws.Range("J1").Value = "ExcludeFields"

# K1: HttpStatusCode -> StatusCode
$ws.Range("K1").Value = "StatusCode"

# Normalize the status code sample values (still 200, just re-written)
$ws.Range("K2").Value = 200
$ws.Range("K3").Value = 200

# The new, slightly longer "ExcludeFields" header needs a bit more room -
# widen column J (10th column) to fit it.
$ws.Columns.Item(10).ColumnWidth = 17.6

# Leave the active cell on L1, just past the last used header column.
$ws.Range("L1").Select()
